$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the "fig 2" row entirely (row 4: photo sieve fig / show sed
#    measurements). Deleting the whole row shifts everything below it up by
#    one, which re-numbers "fig 3".."fig 12" down by one row automatically.
# ---------------------------------------------------------------------------
$ws.Rows("4").Delete()

# ---------------------------------------------------------------------------
# 2) Row 3 (fig 1): retitle the map figure and add a new "Notes" entry.
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "map of lc 1 and 3, study locations"
$ws.Range("E3").Value = "model after the one I did for other paper"

# ---------------------------------------------------------------------------
# 3) Row 4 (fig 3): mark as done.
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = "X"

# ---------------------------------------------------------------------------
# 4) Row 12 (fig 11): add a note, highlighted in red (bold for the label
#    column, regular red for the rest of the row's used cells).
# ---------------------------------------------------------------------------
$ws.Range("E12").Value = "no control / not related… but there is no 40m bedrock exposure so maybe ask sophia to do that?"
$ws.Range("A12").Font.Bold = $true
$ws.Range("A12").Font.Color = 255
$ws.Range("B12").Font.Color = 255
$ws.Range("E12").Font.Color = 255

# ---------------------------------------------------------------------------
# 5) Row 14 (distance up channel vs curvature): mark done and add a red note.
# ---------------------------------------------------------------------------
$ws.Range("D14").Value = "X"
$ws.Range("E14").Value = "no control / not related… maybe shows only bed thickness controls curve"
$ws.Range("A14").Font.Bold = $true
$ws.Range("A14").Font.Color = 255
$ws.Range("B14").Font.Color = 255
$ws.Range("D14").Font.Color = 255
$ws.Range("E14").Font.Color = 255

# ---------------------------------------------------------------------------
# 6) Row 18 (bed thickness distributions ...): mark done.
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = "X"

# ---------------------------------------------------------------------------
# 7) Row 21 (SOIL PIT STUFF): add a red note.
# ---------------------------------------------------------------------------
$ws.Range("E21").Value = "did not work"
$ws.Range("A21").Font.Bold = $true
$ws.Range("A21").Font.Color = 255
$ws.Range("B21").Font.Color = 255
$ws.Range("E21").Font.Color = 255

# ---------------------------------------------------------------------------
# 8) Restore the cursor / selection to B4, matching the author's final state.
# ---------------------------------------------------------------------------
$ws.Range("B4").Select() | Out-Null
